$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New business-day "Serie" labels (August 2021) appended below the existing
# data. They must land in the shared-string table as plain text, exactly
# like the pre-existing "dd-mm-yyyy" labels above them -- so each literal is
# entered with a leading quote to stop Excel's automatic date detection,
# then the quote-prefix formatting is cleared off the whole block again.
$dates = @(
    "02-08-2021","03-08-2021","04-08-2021","05-08-2021","06-08-2021",
    "09-08-2021","10-08-2021","11-08-2021","12-08-2021","13-08-2021",
    "16-08-2021","17-08-2021","18-08-2021","19-08-2021","20-08-2021",
    "23-08-2021","24-08-2021","25-08-2021","26-08-2021","27-08-2021",
    "30-08-2021","31-08-2021"
)

# Data matrix for columns B..J for each new row (Total, 3m, 6m, 9m, 12m, 18m, 2a, 5a, 10a+)
$values = @(
    @(19,0,0,0,0,0,0,16,3),
    @(114,0,48,0,0,27,26,0,14),
    @(524,0,48,135,61,147,16,94,23),
    @(275,0,100,0,94,26,43,5,8),
    @(555,153,128,32,24,0,145,27,46),
    @(422,96,0,0,164,79,55,21,7),
    @(382,0,48,125,87,0,62,29,31),
    @(819,450,0,52,146,0,76,44,51),
    @(417,0,134,103,77,16,79,5,4),
    @(223,0,48,53,38,0,42,37,5),
    @(354,0,0,32,69,0,213,37,2),
    @(80,0,46,0,0,0,14,0,21),
    @(492,0,46,32,118,26,123,55,92),
    @(246,0,48,32,48,0,28,52,38),
    @(323,0,80,84,0,0,48,84,26),
    @(273,0,0,0,32,10,128,75,28),
    @(450,307,48,0,24,0,26,5,41),
    @(422,0,0,84,116,0,79,116,27),
    @(708,0,0,304,158,106,124,16,0),
    @(443,0,79,0,71,210,84,0,0),
    @(254,0,23,0,0,108,78,25,19),
    @(354,0,0,32,142,0,108,19,53)
)

$startRow = 147
$endRow = $startRow + $dates.Count - 1

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = "'" + $dates[$i]
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($row, 2 + $j).Value = $rowVals[$j]
    }
}

# Drop the quote-prefix formatting picked up above so the new "Serie" cells
# end up stored the same plain way as the rest of column A.
$ws.Range("A" + $startRow + ":A" + $endRow).ClearFormats()
